$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    @("MCH189-1", "", "PAPERS, CORRESPONDENCE", "", "Series", "1 Box", "LOCATION: 23P | GRAP COUNT NUMER: NONE", ""),
    @("MCH189-2", "", "CORRESPONDENCE, ANC STATEMENTS 1979, AMANDLA GROUP IN THE USSR, SOVIET POLICY ON SA, ANC OFFICE MOSCOW 1987-93, STYDENT AFFAIRS USSR", "1979", "Series", "1 Box", "LOCATION: 23P | GRAP COUNT NUMER: NONE", ""),
    @("MCH189-3", "", "", "", "Series", "1 Box", "LOCATION: 23P | GRAP COUNT NUMER: NONE", ""),
    @("MCH189-4", "", "CORRESPONDENCE", "", "Series", "1 Box", "LOCATION: 23P | GRAP COUNT NUMER: NONE", ""),
    @("MCH189-5", "", "CORRESPONDENCE", "", "Series", "1 Box", "LOCATION: 23P | GRAP COUNT NUMER: NONE", ""),
    @("MCH189-6", "", "PHOTO ALBUM, AMANDLA GROUP, MANDELA BADGES, POWER TO ... BADGES, MANDELA STAMPS, LETTER TO OR  AND DUMA 1963, PUBLICATIONS", "1963", "Series", "1 Box", "LOCATION: 23P | GRAP COUNT NUMER: NONE", ""),
    @("MCH189-7", "", "AFRO- ASIAN STUDENT COMMITTEE CORRESPONDNCE, ANC STATEMENTS", "", "Series", "1 Box", "LOCATION: 23P | GRAP COUNT NUMER: NONE", "")
)

$r = 2
foreach ($row in $data) {
    $ws.Cells.Item($r, 1).Value = $row[0]
    $ws.Cells.Item($r, 3).Value = $row[2]
    if ($row[3] -ne "") {
        $ws.Cells.Item($r, 4).NumberFormat = "@"
        $ws.Cells.Item($r, 4).Value = $row[3]
    }
    $ws.Cells.Item($r, 5).Value = $row[4]
    $ws.Cells.Item($r, 6).Value = $row[5]
    $ws.Cells.Item($r, 7).Value = $row[6]
    $r = $r + 1
}
